$wb = $excel.ActiveWorkbook

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 5051.1763
$ws_ARM.Range("I32").Value = 3505.625
$ws_ARM.Range("J32").Value = 29780
$ws_ARM.Range("K32").Value = 3505.625
$ws_ARM.Range("L32").Value = 29780
$ws_ARM.Range("M32").Value = -3218.625
$ws_ARM.Range("N32").Value = -30354
$ws_ARM.Range("H63").Value = 3000
$ws_ARM.Range("I63").Value = 3000
$ws_ARM.Range("J63").Value = 3000
$ws_ARM.Range("K63").Value = 3000
$ws_ARM.Range("L63").Value = 3000
$ws_ARM.Range("M63").Value = -2314
$ws_ARM.Range("N63").Value = -4372
$ws_ARM.Range("H66").Value = 3000
$ws_ARM.Range("I66").Value = 3000
$ws_ARM.Range("J66").Value = 3000
$ws_ARM.Range("K66").Value = 15000
$ws_ARM.Range("L66").Value = 15000
$ws_ARM.Range("M66").Value = -11568
$ws_ARM.Range("N66").Value = -21864
$ws_ARM.Range("H74").Value = 3448
$ws_ARM.Range("I74").Value = 3514.6978
$ws_ARM.Range("J74").Value = 2014
$ws_ARM.Range("K74").Value = 3514.6978
$ws_ARM.Range("L74").Value = 2014
$ws_ARM.Range("M74").Value = -2640.6978
$ws_ARM.Range("N74").Value = -3762
$ws_ARM.Range("H77").Value = 3448
$ws_ARM.Range("I77").Value = 3514.6978
$ws_ARM.Range("J77").Value = 2014
$ws_ARM.Range("K77").Value = 17573.489
$ws_ARM.Range("L77").Value = 10070
$ws_ARM.Range("M77").Value = -13205.489
$ws_ARM.Range("N77").Value = -18806
$ws_ARM.Range("H132").Value = 2738.3547
$ws_ARM.Range("I132").Value = 1274.25
$ws_ARM.Range("J132").Value = 3663.0527
$ws_ARM.Range("K132").Value = 3822.75
$ws_ARM.Range("L132").Value = 10989.1581
$ws_ARM.Range("M132").Value = -1292.75
$ws_ARM.Range("N132").Value = -16049.1581
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H20").Value = 11001.333
$ws_BSM.Range("I20").Value = 2916.3333
$ws_BSM.Range("K20").Value = 2916.3333
$ws_BSM.Range("M20").Value = -2669.3333
$ws_BSM.Range("H80").Value = 2317.1714
$ws_BSM.Range("I80").Value = 633.8
$ws_BSM.Range("J80").Value = 3579.7
$ws_BSM.Range("K80").Value = 633.8
$ws_BSM.Range("L80").Value = 3579.7
$ws_BSM.Range("M80").Value = 364.2
$ws_BSM.Range("N80").Value = -5575.7
$ws_BSM.Range("H83").Value = 2317.1714
$ws_BSM.Range("I83").Value = 633.8
$ws_BSM.Range("J83").Value = 3579.7
$ws_BSM.Range("K83").Value = 3169
$ws_BSM.Range("L83").Value = 17898.5
$ws_BSM.Range("M83").Value = 1823
$ws_BSM.Range("N83").Value = -27882.5
$ws_BSM.Range("H94").Value = 761.25
$ws_BSM.Range("I94").Value = 875.625
$ws_BSM.Range("J94").Value = 532.5
$ws_BSM.Range("K94").Value = 875.625
$ws_BSM.Range("L94").Value = 532.5
$ws_BSM.Range("M94").Value = -424.625
$ws_BSM.Range("N94").Value = -1434.5
$ws_BSM.Range("H99").Value = 809.4737
$ws_BSM.Range("I99").Value = 723.75
$ws_BSM.Range("K99").Value = 723.75
$ws_BSM.Range("M99").Value = 774.25
$ws_BSM.Range("H105").Value = 3171.4285
$ws_BSM.Range("I105").Value = 4040
$ws_BSM.Range("K105").Value = 4040
$ws_BSM.Range("M105").Value = -2293
$ws_BSM.Range("H134").Value = 2384.6667
$ws_BSM.Range("I134").Value = 1447.6364
$ws_BSM.Range("K134").Value = 4342.9092
$ws_BSM.Range("M134").Value = -1807.9092
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H58").Value = 1209.6666
$ws_CRP.Range("I58").Value = 1047.2941
$ws_CRP.Range("J58").Value = 2589.8333
$ws_CRP.Range("K58").Value = 1047.2941
$ws_CRP.Range("L58").Value = 2589.8333
$ws_CRP.Range("M58").Value = -844.2941000000001
$ws_CRP.Range("N58").Value = -2995.8333
$ws_CRP.Range("H62").Value = 5093.9165
$ws_CRP.Range("I62").Value = 4029.5
$ws_CRP.Range("J62").Value = 6158.3335
$ws_CRP.Range("K62").Value = 4029.5
$ws_CRP.Range("L62").Value = 6158.3335
$ws_CRP.Range("M62").Value = -3405.5
$ws_CRP.Range("N62").Value = -7406.3335
$ws_CRP.Range("H65").Value = 5093.9165
$ws_CRP.Range("I65").Value = 4029.5
$ws_CRP.Range("J65").Value = 6158.3335
$ws_CRP.Range("K65").Value = 20147.5
$ws_CRP.Range("L65").Value = 30791.6675
$ws_CRP.Range("M65").Value = -17027.5
$ws_CRP.Range("N65").Value = -37031.6675
$ws_CRP.Range("H136").Value = 1209.6666
$ws_CRP.Range("I136").Value = 1047.2941
$ws_CRP.Range("J136").Value = 2589.8333
$ws_CRP.Range("K136").Value = 3141.8823
$ws_CRP.Range("L136").Value = 7769.499899999999
$ws_CRP.Range("M136").Value = -591.8823000000002
$ws_CRP.Range("N136").Value = -12869.4999
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H69").Value = 2420.875
$ws_CUL.Range("I69").Value = 3006
$ws_CUL.Range("J69").Value = 2225.8333
$ws_CUL.Range("K69").Value = 9018
$ws_CUL.Range("L69").Value = 6677.499899999999
$ws_CUL.Range("M69").Value = -8207
$ws_CUL.Range("N69").Value = -8299.499899999999
$ws_CUL.Range("H72").Value = 2420.875
$ws_CUL.Range("I72").Value = 3006
$ws_CUL.Range("J72").Value = 2225.8333
$ws_CUL.Range("K72").Value = 27054
$ws_CUL.Range("L72").Value = 20032.4997
$ws_CUL.Range("M72").Value = -22998
$ws_CUL.Range("N72").Value = -28144.4997
$ws_CUL.Range("H131").Value = 6832.222
$ws_CUL.Range("J131").Value = 9933.333000000001
$ws_CUL.Range("L131").Value = 29799.999
$ws_CUL.Range("N131").Value = -39879.999
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H97").Value = 1404.4546
$ws_GSM.Range("I97").Value = 1534.9
$ws_GSM.Range("J97").Value = 100
$ws_GSM.Range("K97").Value = 1534.9
$ws_GSM.Range("L97").Value = 100
$ws_GSM.Range("M97").Value = -1038.9
$ws_GSM.Range("N97").Value = -1092
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H22").Value = 463.64285
$ws_LTW.Range("I22").Value = 501
$ws_LTW.Range("J22").Value = 426.2857
$ws_LTW.Range("K22").Value = 501
$ws_LTW.Range("L22").Value = 426.2857
$ws_LTW.Range("M22").Value = -206
$ws_LTW.Range("N22").Value = -1016.2857
$ws_LTW.Range("H27").Value = 463.64285
$ws_LTW.Range("I27").Value = 501
$ws_LTW.Range("J27").Value = 426.2857
$ws_LTW.Range("K27").Value = 501
$ws_LTW.Range("L27").Value = 426.2857
$ws_LTW.Range("M27").Value = -394
$ws_LTW.Range("N27").Value = -640.2857
$ws_LTW.Range("H55").Value = 891.1429000000001
$ws_LTW.Range("I55").Value = 516
$ws_LTW.Range("J55").Value = 1172.5
$ws_LTW.Range("K55").Value = 516
$ws_LTW.Range("L55").Value = 1172.5
$ws_LTW.Range("M55").Value = -343
$ws_LTW.Range("N55").Value = -1518.5
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H81").Value = 1666.6428
$ws_WVR.Range("I81").Value = 987.25
$ws_WVR.Range("J81").Value = 1938.4
$ws_WVR.Range("K81").Value = 1974.5
$ws_WVR.Range("L81").Value = 3876.8
$ws_WVR.Range("M81").Value = -913.5
$ws_WVR.Range("N81").Value = -5998.8
$ws_WVR.Range("H84").Value = 1666.6428
$ws_WVR.Range("I84").Value = 987.25
$ws_WVR.Range("J84").Value = 1938.4
$ws_WVR.Range("K84").Value = 9872.5
$ws_WVR.Range("L84").Value = 19384
$ws_WVR.Range("M84").Value = -4568.5
$ws_WVR.Range("N84").Value = -29992
